# Set "想去人数" (column F, rows 2-9) to 0 on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 9; $row++) {
        $ws.Cells.Item($row, 6).Value = 0
    }
}
